$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.771.67'
$ws.Range('E2').Value = '  +0.43%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.824.93'
$ws.Range('E3').Value = '  +1.66%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '349.91'
$ws.Range('E5').Value = '  -0.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '112.78'
$ws.Range('E6').Value = '  +4.02%  '
$ws.Range('E7').Value = '  +1.34%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +3.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.13'
$ws.Range('E10').Value = '  +1.09%  '
$ws.Range('E11').Value = '  -0.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0846'
$ws.Range('E12').Value = '  +1.35%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.08'
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('E14').Value = '  +1.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.273.06'
$ws.Range('E15').Value = '  +1.82%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.981'
$ws.Range('E16').Value = '  +6.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.834.49'
$ws.Range('E17').Value = '  +2.53%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.790.47'
$ws.Range('E18').Value = '  +0.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.42'
$ws.Range('E19').Value = '  +10.00%  '
$ws.Range('E20').Value = '  -0.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.40'
$ws.Range('E21').Value = '  +2.15%  '
$ws.Range('E22').Value = '  +1.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.39'
$ws.Range('E23').Value = '  +0.65%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '268.82'
$ws.Range('E24').Value = '  +0.84%  '
$ws.Range('E25').Value = '  +1.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.23'
$ws.Range('E26').Value = '  +0.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('E28').Value = '  +1.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '39.02'
$ws.Range('E29').Value = '  +7.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.51'
$ws.Range('E30').Value = '  +2.76%  '
$ws.Range('E31').Value = '  +1.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.32'
$ws.Range('E32').Value = '  +1.62%  '
$ws.Range('E33').Value = '  +1.71%  '
$ws.Range('E34').Value = '  +7.99%  '
$ws.Range('E35').Value = '  +2.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0449'
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.84'
$ws.Range('E38').Value = '  +1.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.21'
$ws.Range('E39').Value = '  +2.05%  '
$ws.Range('E40').Value = '  +2.15%  '
$ws.Range('E41').Value = '  +1.31%  '
$ws.Range('E42').Value = '  -1.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '122.91'
$ws.Range('E43').Value = '  +2.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '22.07'
$ws.Range('E45').Value = '  +0.33%  '
$ws.Range('E46').Value = '  +9.13%  '
$ws.Range('E47').Value = '  +7.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.168.67'
$ws.Range('E48').Value = '  +2.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.249'
$ws.Range('E49').Value = '  +22.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.941'
$ws.Range('E50').Value = '  +4.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.52'
$ws.Range('E51').Value = '  +1.66%  '
